$wb = $excel.ActiveWorkbook

# Sheet: 展览 (Exhibition)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 7714
$ws1.Range("F3").Value = 72
$ws1.Range("F4").Value = 222
$ws1.Range("F5").Value = 48
$ws1.Range("F6").Value = 422
$ws1.Range("F7").Value = 1151
$ws1.Range("F8").Value = 204

# Sheet: 演出 (Performance)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 13

# Sheet: 全部类型 (All Types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 7714
$ws4.Range("F3").Value = 72
$ws4.Range("F4").Value = 222
$ws4.Range("F5").Value = 48
$ws4.Range("F6").Value = 422
$ws4.Range("F7").Value = 1151
$ws4.Range("F8").Value = 204
$ws4.Range("F9").Value = 13
